# "some changes to the prnting" - fill in the previously-empty minPts=8 row
# on the "All" and "Without Outliers" sheets (the underlying stats run was
# re-extended to include an extra minPts value), which in turn makes the
# "Ratios" sheet's row 8 (= 'Without Outliers'/All) resolve to real numbers
# instead of #DIV/0! now that the denominators are populated. Also a couple
# of previously-computed error cells changed from #VALUE!/#NUM! to #N/A
# (the stats engine now reports "not available" rather than a value/num
# error for those particular cells), and the active selection moved down
# one row (from the row 7 stats row to the new row 8 stats row).

$wb = $excel.ActiveWorkbook

$all = $wb.Worksheets.Item("All")
$wo  = $wb.Worksheets.Item("Without Outliers")
$ratios = $wb.Worksheets.Item("Ratios")

# ---------------------------------------------------------------------
# 1. New minPts=8 row of statistics on "All"
# ---------------------------------------------------------------------
$all.Range("B8").Value = 1.80446
$all.Range("C8").Value = 0.333857828174409
$all.Range("D8").Value = "#VALUE!"
$all.Range("E8").Value = 1.65575
$all.Range("F8").Value = 1.485415
$all.Range("G8").Value = 2.04915
$all.Range("H8").Value = 0.3343831483
$all.Range("I8").Value = 0.578258720902677
$all.Range("J8").Value = "#NUM!"
$all.Range("K8").Value = 1.08072268797666
$all.Range("L8").Value = 1.12747
$all.Range("M8").Value = 1.31508
$all.Range("N8").Value = 2.44255
$all.Range("O8").Value = 5.41338
$all.Range("P8").Value = 3

# ---------------------------------------------------------------------
# 2. New minPts=8 row of statistics on "Without Outliers"
# ---------------------------------------------------------------------
$wo.Range("B8").Value = 1.485415
$wo.Range("C8").Value = 0.170335
$wo.Range("D8").Value = "#VALUE!"
$wo.Range("E8").Value = 1.485415
$wo.Range("F8").Value = 1.4002475
$wo.Range("G8").Value = 1.5705825
$wo.Range("H8").Value = 0.05802802445
$wo.Range("I8").Value = 0.240890067146821
$wo.Range("J8").Value = "#NUM!"
$wo.Range("K8").Value = "#NUM!"
$wo.Range("L8").Value = 0.34067
$wo.Range("M8").Value = 1.31508
$wo.Range("N8").Value = 1.65575
$wo.Range("O8").Value = 2.97083
$wo.Range("P8").Value = 2

# ---------------------------------------------------------------------
# 3. A handful of previously-cached error results switched from
#    #VALUE!/#NUM! to #N/A on "All" and "Without Outliers" (rows 3-7).
#    These cells hold literal error values, not live formulas.
# ---------------------------------------------------------------------
$all.Range("D3").Value = "#N/A"
$all.Range("D4").Value = "#N/A"
$all.Range("D5").Value = "#N/A"
$all.Range("D6").Value = "#N/A"
$all.Range("D7").Value = "#N/A"
$all.Range("J3").Value = "#N/A"
$all.Range("J4").Value = "#N/A"

$wo.Range("D3").Value = "#N/A"
$wo.Range("D4").Value = "#N/A"
$wo.Range("D5").Value = "#N/A"
$wo.Range("D6").Value = "#N/A"
$wo.Range("D7").Value = "#N/A"
$wo.Range("J3").Value = "#N/A"
$wo.Range("J4").Value = "#N/A"
$wo.Range("J6").Value = "#N/A"
$wo.Range("J7").Value = "#N/A"
$wo.Range("K3").Value = "#N/A"
$wo.Range("K4").Value = "#N/A"

# ---------------------------------------------------------------------
# 4. Move the on-sheet selection down to the newly-filled row (B8:P8),
#    matching where the user's cursor ended up. "Without Outliers" is
#    the active tab, so update it last without disturbing that.
# ---------------------------------------------------------------------
$ratios.Range("B8:P8").Select()
$all.Range("B8:P8").Select()
$wo.Range("B8:P8").Select()
